$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so we can write the updated figures,
# then restore protection afterwards.
$ws.Unprotect()

# Update the confidential disclosure text: model date 2021-05-13 -> 2021-05-14
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-14 for illustrative purposes only and are subject to change."

# Updated holdings weight / percent-change figures
$ws.Range("D2").Value = 0.849441666244509
$ws.Range("E2").Value = 0.01674121405750806
$ws.Range("D3").Value = 0.1505583337554911
$ws.Range("E3").Value = 0.018073447413959
$ws.Range("E4").Value = 0.01694179289182873

# Row 7 grows to accommodate the two-line disclosure; let Excel auto-size it
# back rather than leaving a stale explicit height.
$ws.Rows(7).AutoFit()

$ws.Protect()
